# Apply the edits described in the commit:
# "Decreased probability of fetal death before 4 weeks from conception"
# Updates the B2:B5 probabilities (and dependent D formulas) on sheet "Phase1",
# and the C2:C5 probabilities (and dependent E formulas) on sheet "Phase2".
# Also updates the active selection on each sheet to reflect the edited ranges.

$wb = $excel.ActiveWorkbook

# --- Phase1 sheet ---
$ws1 = $wb.Worksheets.Item("Phase1")

$ws1.Range("B2").Value = 0.1
$ws1.Range("B3").Value = 0.1
$ws1.Range("B4").Value = 0.05
$ws1.Range("B5").Value = 0.05

$ws1.Select()
$ws1.Range("B2:B5").Select()

# --- Phase2 sheet ---
$ws2 = $wb.Worksheets.Item("Phase2")

$ws2.Range("C2").Value = 0.1
$ws2.Range("C3").Value = 0.1
$ws2.Range("C4").Value = 0.05
$ws2.Range("C5").Value = 0.05

$ws2.Select()
$ws2.Range("C2:C5").Select()

$wb.Save()
